$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C) for all data rows (2-89)
# from 2023-09-15 (45184) to 2023-09-17 (45186)
for ($r = 2; $r -le 89; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45184) {
        $cell.Value = 45186
    }
}

# Columns S, T, V, W, X, Y contain HYPERLINK formulas; add a friendly-name
# second argument (the "Beteckning" in column A) to each formula, rows 2-12.
for ($r = 2; $r -le 12; $r++) {
    $beteckning = $ws.Cells.Item($r, 1).Value2

    $artfyndUrl = "https://klasma.github.io/Logging_UPPLANDS-BRO/artfynd/$beteckning.xlsx"
    $ws.Range("S$r").Formula = "=HYPERLINK(`"$artfyndUrl`", `"$beteckning`")"

    $kartorUrl = "https://klasma.github.io/Logging_UPPLANDS-BRO/kartor/$beteckning.png"
    $ws.Range("T$r").Formula = "=HYPERLINK(`"$kartorUrl`", `"$beteckning`")"

    $klagomalUrl = "https://klasma.github.io/Logging_UPPLANDS-BRO/klagomål/$beteckning.docx"
    $ws.Range("V$r").Formula = "=HYPERLINK(`"$klagomalUrl`", `"$beteckning`")"

    $klagomalsmailUrl = "https://klasma.github.io/Logging_UPPLANDS-BRO/klagomålsmail/$beteckning.docx"
    $ws.Range("W$r").Formula = "=HYPERLINK(`"$klagomalsmailUrl`", `"$beteckning`")"

    $tillsynUrl = "https://klasma.github.io/Logging_UPPLANDS-BRO/tillsyn/$beteckning.docx"
    $ws.Range("X$r").Formula = "=HYPERLINK(`"$tillsynUrl`", `"$beteckning`")"

    $tillsynsmailUrl = "https://klasma.github.io/Logging_UPPLANDS-BRO/tillsynsmail/$beteckning.docx"
    $ws.Range("Y$r").Formula = "=HYPERLINK(`"$tillsynsmailUrl`", `"$beteckning`")"
}
